# Generate Report for Handoff
# Updates the localization-status workbook: the markdown file that was
# handed off gets a new GUID/timestamp, and two new dependency PNG files
# show up as additional rows on every sheet.

$wb = $excel.ActiveWorkbook

$oldGuid = "2b2212d2-b2b6-427b-8c63-d0e49f0dada0"
$newGuid = "6836bd44-93f7-4dc5-8afe-c5db7d0c1d8b"

$oldZhXlf = "$oldGuid.284503fe56eeb9e43bf9711ce3cced7004368416.zh-cn.xlf"
$newZhXlf = "$newGuid.76e7c0c1eb084f829b453c5184f7be8dffcf9d2c.zh-cn.xlf"
$oldDeXlf = "$oldGuid.284503fe56eeb9e43bf9711ce3cced7004368416.de-de.xlf"
$newDeXlf = "$newGuid.76e7c0c1eb084f829b453c5184f7be8dffcf9d2c.de-de.xlf"

$png1 = "89ad060a-b4ea-4333-affe-d47088e31631.png"
$png2 = "d47065c8-712b-4bd0-a75e-228caba8ab79.png"
$png1Zh = "ce7b7270b8ca45039a6e3816b4554a418b6191ea.png"
$png2Zh = "cbe8bc1048af02662d6a843c041c160c5bad03a9.png"

$newHandoffDate = "2016-03-21 19:00:23"
$newHandoffDatetimeZh = "2016-03-21 19:00:19"
$epoch = "0001-01-01 00:00:00"

$baseMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/21aa9b6a144dae10060f2d090b3f68afe6f0a3ae/e2e"
$baseZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/014dd36656e1f010c88c8e92057b1d5127c9f50a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$baseDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/437a8f19f9b731c296c26f3a083782e1dfdf38ad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

function Set-LinkCell($ws, $cell, $text, $url) {
    $ws.Hyperlinks.Add($cell, $url, "", "", $text) | Out-Null
}

function Set-DateTextCell($cell, $text) {
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $cell.Value = $text
}

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Hyperlinks.Add always *appends* rather than replacing any existing
# hyperlink on the same cell, and Range.Hyperlinks.Delete() clears every
# hyperlink on the whole sheet (not just the range). So wipe all of this
# sheet's hyperlinks once up front, then (re)add every one of them below,
# in the same left-to-right / top-to-bottom order they should appear in.
$wsOverview.Range("A1").Hyperlinks.Delete()

# Update existing handoff row (row 2) for the markdown file.
Set-LinkCell $wsOverview $wsOverview.Range("A2") "$newGuid.md" "$baseMdUrl/$newGuid.md"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
Set-DateTextCell $wsOverview.Range("D2") $newHandoffDate

# New rows for the two dependency png files.
Set-LinkCell $wsOverview $wsOverview.Range("A3") $png1 "$baseMdUrl/$png1"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
Set-DateTextCell $wsOverview.Range("D3") $newHandoffDate

Set-LinkCell $wsOverview $wsOverview.Range("A4") $png2 "$baseMdUrl/$png2"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
Set-DateTextCell $wsOverview.Range("D4") $newHandoffDate

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A1").Hyperlinks.Delete()

Set-LinkCell $wsZh $wsZh.Range("A2") "$newGuid.md" "$baseMdUrl/$newGuid.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Ready for handoff"
Set-LinkCell $wsZh $wsZh.Range("D2") $newZhXlf "$baseZhUrl/$newZhXlf"
Set-DateTextCell $wsZh.Range("E2") $newHandoffDatetimeZh
Set-DateTextCell $wsZh.Range("H2") $epoch
$wsZh.Range("J2").Value = "Include"

Set-LinkCell $wsZh $wsZh.Range("A3") $png1 "$baseMdUrl/$png1"
$wsZh.Range("B3").Value = ".png"
$wsZh.Range("C3").Value = "Ready for handoff"
Set-LinkCell $wsZh $wsZh.Range("D3") $png1Zh "$baseZhUrl/$png1Zh"
Set-DateTextCell $wsZh.Range("E3") $newHandoffDatetimeZh
Set-DateTextCell $wsZh.Range("H3") $epoch
$wsZh.Range("J3").Value = "IsDependency"
$wsZh.Range("K3").Value = "e2e\$newGuid.md"

Set-LinkCell $wsZh $wsZh.Range("A4") $png2 "$baseMdUrl/$png2"
$wsZh.Range("B4").Value = ".png"
$wsZh.Range("C4").Value = "Ready for handoff"
Set-LinkCell $wsZh $wsZh.Range("D4") $png2Zh "$baseZhUrl/$png2Zh"
Set-DateTextCell $wsZh.Range("E4") $newHandoffDatetimeZh
Set-DateTextCell $wsZh.Range("H4") $epoch
$wsZh.Range("J4").Value = "IsDependency"
$wsZh.Range("K4").Value = "e2e\$newGuid.md"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A1").Hyperlinks.Delete()

Set-LinkCell $wsDe $wsDe.Range("A2") "$newGuid.md" "$baseMdUrl/$newGuid.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Ready for handoff"
Set-LinkCell $wsDe $wsDe.Range("D2") $newDeXlf "$baseDeUrl/$newDeXlf"
Set-DateTextCell $wsDe.Range("E2") $newHandoffDate
Set-DateTextCell $wsDe.Range("H2") $epoch
$wsDe.Range("J2").Value = "Include"

Set-LinkCell $wsDe $wsDe.Range("A3") $png1 "$baseMdUrl/$png1"
$wsDe.Range("B3").Value = ".png"
$wsDe.Range("C3").Value = "Ready for handoff"
Set-LinkCell $wsDe $wsDe.Range("D3") $png1Zh "$baseDeUrl/$png1Zh"
Set-DateTextCell $wsDe.Range("E3") $newHandoffDate
Set-DateTextCell $wsDe.Range("H3") $epoch
$wsDe.Range("J3").Value = "IsDependency"
$wsDe.Range("K3").Value = "e2e\$newGuid.md"

Set-LinkCell $wsDe $wsDe.Range("A4") $png2 "$baseMdUrl/$png2"
$wsDe.Range("B4").Value = ".png"
$wsDe.Range("C4").Value = "Ready for handoff"
Set-LinkCell $wsDe $wsDe.Range("D4") $png2Zh "$baseDeUrl/$png2Zh"
Set-DateTextCell $wsDe.Range("E4") $newHandoffDate
Set-DateTextCell $wsDe.Range("H4") $epoch
$wsDe.Range("J4").Value = "IsDependency"
$wsDe.Range("K4").Value = "e2e\$newGuid.md"

Write-Host "Report generated for handoff."
